# Updated symbol list on Wed Dec 21 21:36:17 UTC 2022 with GitHub Actions
# Refresh the "Price" (column D) values for the listed coins.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = "247.03"
    3  = "22.37"
    4  = "5.236"
    5  = "0.05694"
    6  = "3.414"
    7  = "6.305"
    8  = "0.8056"
    9  = "0.8596"
    11 = "0.07372"
    12 = "0.03037"
    13 = "0.03096"
    14 = "0.09386"
    15 = "3.867"
    16 = "0.001574"
    17 = "0.04768"
    19 = "0.006377"
    21 = "0.0009964"
    22 = "0.0001500"
    23 = "3.691"
    24 = "2.194"
    40 = "0.03918"
    42 = "0.1065"
    44 = "0.008429"
    45 = "0.00005593"
    47 = "0.4499"
    48 = "0.1851"
    49 = "0.00002100"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    # Force the numeric-looking price to be stored as text (matching the
    # original "Price" column, which is text) without leaving a quote-prefix
    # or a lingering custom number-format style behind.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    $cell.ClearFormats()
}
